# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the e0d7219b-...-37081b86aeef entry on each sheet, as produced by a
# fresh handback-status report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     e0d7219b... row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 12:56:25"

# --- zh-cn sheet: bump "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (K) for the e0d7219b... row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-02 12:56:20"
$wsZhCn.Range("K3").Value = "2016-09-02 12:56:38"

# --- de-de sheet: bump "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (K) for the e0d7219b... row (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-02 12:56:25"
$wsDeDe.Range("K3").Value = "2016-09-02 12:56:44"
